# Update of files: refresh "Förändrad" end-date and re-sync the
# Beteckning/Datum/Markägare/Area rows for sheet "Avverkningsanmälningar".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") advances by one day (2026-02-06 -> 2026-02-07 / 46059 -> 46060)
# for every data row.
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 3).Value = 46060
}

# Row 5: Beteckning/Datum/Markägare/Area now matches former row 8 ("A 42462-2025")
$ws.Cells.Item(5, 1).Value = "A 42462-2025"
$ws.Cells.Item(5, 2).Value = 45905.45730324074
$ws.Cells.Item(5, 6).ClearContents() | Out-Null
$ws.Cells.Item(5, 7).Value = 2.3

# Row 6: Beteckning/Datum/Markägare/Area now matches former row 38 ("A 1103-2025")
$ws.Cells.Item(6, 1).Value = "A 1103-2025"
$ws.Cells.Item(6, 2).Value = 45666.0
$ws.Cells.Item(6, 6).ClearContents() | Out-Null
$ws.Cells.Item(6, 7).Value = 1.6

# Row 7: Beteckning/Datum/Markägare/Area now matches former row 31 ("A 21141-2023")
$ws.Cells.Item(7, 1).Value = "A 21141-2023"
$ws.Cells.Item(7, 2).Value = 45062.0
$ws.Cells.Item(7, 6).ClearContents() | Out-Null
$ws.Cells.Item(7, 7).Value = 3.4

# Row 8: Beteckning/Datum/Markägare/Area now matches former row 12 ("A 45167-2025")
$ws.Cells.Item(8, 1).Value = "A 45167-2025"
$ws.Cells.Item(8, 2).Value = 45919.49364583333
$ws.Cells.Item(8, 6).ClearContents() | Out-Null
$ws.Cells.Item(8, 7).Value = 0.6

# Row 9: Beteckning/Datum/Markägare/Area now matches former row 15 ("A 45088-2025")
$ws.Cells.Item(9, 1).Value = "A 45088-2025"
$ws.Cells.Item(9, 2).Value = 45919.37598379629
$ws.Cells.Item(9, 6).ClearContents() | Out-Null
$ws.Cells.Item(9, 7).Value = 0.9

# Row 10: Beteckning/Datum/Markägare/Area now matches former row 14 ("A 45158-2025")
$ws.Cells.Item(10, 1).Value = "A 45158-2025"
$ws.Cells.Item(10, 2).Value = 45919.482453703706
$ws.Cells.Item(10, 6).ClearContents() | Out-Null
$ws.Cells.Item(10, 7).Value = 2.7

# Row 11: Beteckning/Datum/Markägare/Area now matches former row 19 ("A 1379-2024")
$ws.Cells.Item(11, 1).Value = "A 1379-2024"
$ws.Cells.Item(11, 2).Value = 45303.551932870374
$ws.Cells.Item(11, 6).Value = "Kommuner"
$ws.Cells.Item(11, 7).Value = 1.5

# Row 12: Beteckning/Datum/Markägare/Area now matches former row 21 ("A 61627-2024")
$ws.Cells.Item(12, 1).Value = "A 61627-2024"
$ws.Cells.Item(12, 2).Value = 45646.66263888889
$ws.Cells.Item(12, 6).Value = "Kommuner"
$ws.Cells.Item(12, 7).Value = 2.1

# Row 13: Beteckning/Datum/Markägare/Area now matches former row 27 ("A 34591-2022")
$ws.Cells.Item(13, 1).Value = "A 34591-2022"
$ws.Cells.Item(13, 2).Value = 44795.377858796295
$ws.Cells.Item(13, 6).ClearContents() | Out-Null
$ws.Cells.Item(13, 7).Value = 2.5

# Row 14: Beteckning/Datum/Markägare/Area now matches former row 41 ("A 34508-2025")
$ws.Cells.Item(14, 1).Value = "A 34508-2025"
$ws.Cells.Item(14, 2).Value = 45847.44315972222
$ws.Cells.Item(14, 6).ClearContents() | Out-Null
$ws.Cells.Item(14, 7).Value = 2

# Row 15: Beteckning/Datum/Markägare/Area now matches former row 39 ("A 34394-2025")
$ws.Cells.Item(15, 1).Value = "A 34394-2025"
$ws.Cells.Item(15, 2).Value = 45846.588541666664
$ws.Cells.Item(15, 6).ClearContents() | Out-Null
$ws.Cells.Item(15, 7).Value = 3.9

# Row 16: Beteckning/Datum/Markägare/Area now matches former row 40 ("A 32577-2025")
$ws.Cells.Item(16, 1).Value = "A 32577-2025"
$ws.Cells.Item(16, 2).Value = 45838.0
$ws.Cells.Item(16, 6).ClearContents() | Out-Null
$ws.Cells.Item(16, 7).Value = 5.3

# Row 18: Beteckning/Datum/Markägare/Area now matches former row 5 ("A 34939-2025")
$ws.Cells.Item(18, 1).Value = "A 34939-2025"
$ws.Cells.Item(18, 2).Value = 45849.582291666666
$ws.Cells.Item(18, 6).ClearContents() | Out-Null
$ws.Cells.Item(18, 7).Value = 7.8

# Row 19: Beteckning/Datum/Markägare/Area now matches former row 6 ("A 34963-2025")
$ws.Cells.Item(19, 1).Value = "A 34963-2025"
$ws.Cells.Item(19, 2).Value = 45849.632199074076
$ws.Cells.Item(19, 6).ClearContents() | Out-Null
$ws.Cells.Item(19, 7).Value = 1.1

# Row 20: Beteckning/Datum/Markägare/Area now matches former row 24 ("A 21379-2023")
$ws.Cells.Item(20, 1).Value = "A 21379-2023"
$ws.Cells.Item(20, 2).Value = 45063.34819444444
$ws.Cells.Item(20, 6).ClearContents() | Out-Null
$ws.Cells.Item(20, 7).Value = 5.8

# Row 21: Beteckning/Datum/Markägare/Area now matches former row 26 ("A 57000-2025")
$ws.Cells.Item(21, 1).Value = "A 57000-2025"
$ws.Cells.Item(21, 2).Value = 45977.0
$ws.Cells.Item(21, 6).ClearContents() | Out-Null
$ws.Cells.Item(21, 7).Value = 2.3

# Row 22: Beteckning/Datum/Markägare/Area now matches former row 7 ("A 57655-2025")
$ws.Cells.Item(22, 1).Value = "A 57655-2025"
$ws.Cells.Item(22, 2).Value = 45981.40369212963
$ws.Cells.Item(22, 6).ClearContents() | Out-Null
$ws.Cells.Item(22, 7).Value = 3.4

# Row 23: Beteckning/Datum/Markägare/Area now matches former row 32 ("A 60392-2022")
$ws.Cells.Item(23, 1).Value = "A 60392-2022"
$ws.Cells.Item(23, 2).Value = 44910.0
$ws.Cells.Item(23, 6).ClearContents() | Out-Null
$ws.Cells.Item(23, 7).Value = 4.1

# Row 24: Beteckning/Datum/Markägare/Area now matches former row 10 ("A 59011-2025")
$ws.Cells.Item(24, 1).Value = "A 59011-2025"
$ws.Cells.Item(24, 2).Value = 45987.0
$ws.Cells.Item(24, 6).ClearContents() | Out-Null
$ws.Cells.Item(24, 7).Value = 2.7

# Row 26: Beteckning/Datum/Markägare/Area now matches former row 13 ("A 2434-2026")
$ws.Cells.Item(26, 1).Value = "A 2434-2026"
$ws.Cells.Item(26, 2).Value = 46036.86722222222
$ws.Cells.Item(26, 6).ClearContents() | Out-Null
$ws.Cells.Item(26, 7).Value = 1.2

# Row 27: Beteckning/Datum/Markägare/Area now matches former row 16 ("A 2433-2026")
$ws.Cells.Item(27, 1).Value = "A 2433-2026"
$ws.Cells.Item(27, 2).Value = 46036.866319444445
$ws.Cells.Item(27, 6).ClearContents() | Out-Null
$ws.Cells.Item(27, 7).Value = 2.6

# Row 28: Beteckning/Datum/Markägare/Area now matches former row 34 ("A 60717-2022")
$ws.Cells.Item(28, 1).Value = "A 60717-2022"
$ws.Cells.Item(28, 2).Value = 44912.89078703704
$ws.Cells.Item(28, 6).ClearContents() | Out-Null
$ws.Cells.Item(28, 7).Value = 0.7

# Row 29: Beteckning/Datum/Markägare/Area now matches former row 22 ("A 43229-2024")
$ws.Cells.Item(29, 1).Value = "A 43229-2024"
$ws.Cells.Item(29, 2).Value = 45567.886979166666
$ws.Cells.Item(29, 6).ClearContents() | Out-Null
$ws.Cells.Item(29, 7).Value = 9.7

# Row 30: Beteckning/Datum/Markägare/Area now matches former row 33 ("A 59192-2022")
$ws.Cells.Item(30, 1).Value = "A 59192-2022"
$ws.Cells.Item(30, 2).Value = 44896.0
$ws.Cells.Item(30, 6).ClearContents() | Out-Null
$ws.Cells.Item(30, 7).Value = 1.5

# Row 31: Beteckning/Datum/Markägare/Area now matches former row 11 ("A 7791-2023")
$ws.Cells.Item(31, 1).Value = "A 7791-2023"
$ws.Cells.Item(31, 2).Value = 44973.0
$ws.Cells.Item(31, 6).ClearContents() | Out-Null
$ws.Cells.Item(31, 7).Value = 3.1

# Row 32: Beteckning/Datum/Markägare/Area now matches former row 23 ("A 769-2023")
$ws.Cells.Item(32, 1).Value = "A 769-2023"
$ws.Cells.Item(32, 2).Value = 44931.0
$ws.Cells.Item(32, 6).ClearContents() | Out-Null
$ws.Cells.Item(32, 7).Value = 1.7

# Row 33: Beteckning/Datum/Markägare/Area now matches former row 9 ("A 1486-2023")
$ws.Cells.Item(33, 1).Value = "A 1486-2023"
$ws.Cells.Item(33, 2).Value = 44937.0
$ws.Cells.Item(33, 6).ClearContents() | Out-Null
$ws.Cells.Item(33, 7).Value = 0.2

# Row 34: Beteckning/Datum/Markägare/Area now matches former row 29 ("A 60793-2023")
$ws.Cells.Item(34, 1).Value = "A 60793-2023"
$ws.Cells.Item(34, 2).Value = 45260.6534375
$ws.Cells.Item(34, 6).ClearContents() | Out-Null
$ws.Cells.Item(34, 7).Value = 0.7

# Row 35: Beteckning/Datum/Markägare/Area now matches former row 20 ("A 57410-2022")
$ws.Cells.Item(35, 1).Value = "A 57410-2022"
$ws.Cells.Item(35, 2).Value = 44896.0
$ws.Cells.Item(35, 6).ClearContents() | Out-Null
$ws.Cells.Item(35, 7).Value = 7.5

# Row 36: Beteckning/Datum/Markägare/Area now matches former row 18 ("A 31486-2021")
$ws.Cells.Item(36, 1).Value = "A 31486-2021"
$ws.Cells.Item(36, 2).Value = 44369.43783564815
$ws.Cells.Item(36, 6).ClearContents() | Out-Null
$ws.Cells.Item(36, 7).Value = 4.8

# Row 37: Beteckning/Datum/Markägare/Area now matches former row 28 ("A 63664-2023")
$ws.Cells.Item(37, 1).Value = "A 63664-2023"
$ws.Cells.Item(37, 2).Value = 45275.62074074074
$ws.Cells.Item(37, 6).ClearContents() | Out-Null
$ws.Cells.Item(37, 7).Value = 3.2

# Row 38: Beteckning/Datum/Markägare/Area now matches former row 30 ("A 57955-2024")
$ws.Cells.Item(38, 1).Value = "A 57955-2024"
$ws.Cells.Item(38, 2).Value = 45631.569398148145
$ws.Cells.Item(38, 6).ClearContents() | Out-Null
$ws.Cells.Item(38, 7).Value = 0.9

# Row 39: Beteckning/Datum/Markägare/Area now matches former row 37 ("A 23295-2025")
$ws.Cells.Item(39, 1).Value = "A 23295-2025"
$ws.Cells.Item(39, 2).Value = 45791.59071759259
$ws.Cells.Item(39, 6).ClearContents() | Out-Null
$ws.Cells.Item(39, 7).Value = 2.2

# Row 40: Beteckning/Datum/Markägare/Area now matches former row 36 ("A 23294-2025")
$ws.Cells.Item(40, 1).Value = "A 23294-2025"
$ws.Cells.Item(40, 2).Value = 45791.58967592593
$ws.Cells.Item(40, 6).ClearContents() | Out-Null
$ws.Cells.Item(40, 7).Value = 4.6

# Row 41: Beteckning/Datum/Markägare/Area now matches former row 35 ("A 23301-2025")
$ws.Cells.Item(41, 1).Value = "A 23301-2025"
$ws.Cells.Item(41, 2).Value = 45791.594988425924
$ws.Cells.Item(41, 6).ClearContents() | Out-Null
$ws.Cells.Item(41, 7).Value = 2.9

